$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: DATE COMPLETED changes
$ws.Range("B10").Value = 43912

# Copy the date formatting used by existing date cells (numFmtId 14) onto the new rows
$ws.Range("A9:B9").Copy()
$ws.Range("A11:B13").PasteSpecial(-4122)

# Row 11
$ws.Range("A11").Value = 43908
$ws.Range("B11").Value = 43912
$ws.Range("C11").Value = "leftview binary tree"

# Row 12
$ws.Range("A12").Value = 43908
$ws.Range("B12").Value = 43912
$ws.Range("C12").Value = "level-order traversal"

# Row 13
$ws.Range("A13").Value = 43914
$ws.Range("B13").Value = 43914
$ws.Range("C13").Value = "Activity Selection "
$ws.Range("D13").Value = "https://www.techiedelight.com/activity-selection-problem/"

$ws.Range("B13").Select()
